$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "26.932.19") that look numeric but are
# plain text (thousands-dot-grouped, not valid floating point values). Force
# the cells we are about to rewrite to Text format first so COM does not
# reinterpret values like "9.20" as the number 9.2 and drop the trailing zero.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D23", "D25", "D33", "D39", "D41", "D44", "D47", "D48")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.932.19'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.551.05'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '206.66'
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").Value = '0.486'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("D9").Value = '0.247'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '0.0586'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '1.772.80'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '1.550.28'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").Value = '3.73'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").Value = '26.930.22'
$ws.Range("D17").Value = '61.63'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '217.49'
$ws.Range("E18").Value = '  +1.96%  '
$ws.Range("D19").Value = '0.0₃0696'
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("D23").Value = '9.20'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").Value = '154.12'
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("E32").Value = '  +0.20%  '
$ws.Range("D33").Value = '1.422.77'
$ws.Range("E33").Value = '  +4.34%  '
$ws.Range("E34").Value = '  +4.18%  '
$ws.Range("E35").Value = '  +3.11%  '
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("D39").Value = '0.522'
$ws.Range("E39").Value = '  +1.24%  '
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("D41").Value = '5.74'
$ws.Range("E41").Value = '  +5.19%  '
$ws.Range("E43").Value = '  +4.61%  '
$ws.Range("D44").Value = '0.992'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").Value = '1.685.98'
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").Value = '87.65'
$ws.Range("E48").Value = '  +1.64%  '
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("E50").Value = '  +3.36%  '
$ws.Range("E51").Value = '  +0.29%  '
